$wb = $excel.ActiveWorkbook

# --- Sheet 1: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M4").Value = 95.56
$ws1.Range("M21").Value = "1 de 19"

# --- Sheet 2: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F4").Value = 95.56
$ws2.Range("F21").Value = 463.36

# --- Sheet 3: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Columns.Item(6).ColumnWidth = 25.16666666666667

$ws3.Range("D3").Value = 95.56
$ws3.Range("E3").Value = 15375.9993
$ws3.Range("F3").Value = 0.006176494440350301

$ws3.Range("D4").Value = 70070.33
$ws3.Range("E4").Value = -54598.77070000001
$ws3.Range("F4").Value = 4.528976597724057
